$d = $word.ActiveDocument

# --- Edit 1: "Memoria RAM 2 GB" -> "Memoria RAM 4 GB", split as three runs
#     ("Memoria RAM ", "4", " GB") so the changed digit becomes its own run. ---
$rng = $d.Content
$rng.Find.Execute("Memoria RAM 2 GB", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$start = $rng.Start
$digitStart = $start + 12
$digitEnd = $digitStart + 1
$digitRng = $d.Range($digitStart, $digitEnd)
$digitRng.Text = "4"
# Touch (and revert) formatting on just the new run so it is not silently
# re-coalesced into its neighbours, matching the three-run split in the source.
$touch = $d.Range($digitStart, $digitStart + 1)
$touch.Bold = 1
$touch.Bold = 0

# --- Edit 2: "512Mb" -> "1GB" ---
$rng2 = $d.Content
$rng2.Find.Execute("512Mb", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$mStart = $rng2.Start
$mEnd = $rng2.End
$r2 = $d.Range($mStart, $mEnd)
$r2.Text = "1GB"
$touch2 = $d.Range($mStart, $mStart + 3)
$touch2.Bold = 1
$touch2.Bold = 0

# --- Edit 3: add a new bullet "Sistema operativo Linux" right after the
#     "Sistema operativo Windows 8 o superior " bullet. ---
$rng3 = $d.Content
$rng3.Find.Execute("Sistema operativo Windows 8 o superior", $true, $false, $false, `
                    $false, $false, $true, 1, $false, "", 0) | Out-Null
$winPara = $rng3.Paragraphs(1)
$winPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($winPara.Index + 1)
$newPara.Range.Text = "Sistema operativo Linux"
